# Update Unit Test/XML Templates for onlineResources
#
# Adds new lookup rows describing GMD_ONLINERESOURCE fields (and a
# corrected "contacts{}.postalcode" property) to the Config_AB and
# Config_BC sheets, and leaves the Config_BC sheet as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Config_AB (sheet2)
# ---------------------------------------------------------------------
$wsAB = $wb.Worksheets.Item("Config_AB")

# online_resource_link / online_resource_protocol go in first so the new
# shared strings are created in the same order as the authored workbook.
$wsAB.Range("A76").Value = "online_resource_link"
$wsAB.Range("B76").Value = "online_resource_link"
$wsAB.Range("C76").Value = "GMD_ONLINERESOURCE"

$wsAB.Range("A77").Value = "online_resource_protocol"
$wsAB.Range("B77").Value = "online_resource_protocol"
$wsAB.Range("C77").Value = "GMD_ONLINERESOURCE"

# Contact-info block pasted in directly below (temporarily at rows
# 78-84); it gets pushed down to 81-87 once the description block is
# inserted above it.
$wsAB.Range("A78").Value = "contact_phone"
$wsAB.Range("B78").Value = "contacts{}.phone"
$wsAB.Range("C78").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsAB.Range("A79").Value = "contact_facsimile"
$wsAB.Range("B79").Value = "contacts{}.facsimile"
$wsAB.Range("C79").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsAB.Range("A80").Value = "contact_address"
$wsAB.Range("B80").Value = "contacts{}.address"
$wsAB.Range("C80").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsAB.Range("A81").Value = "contact_city"
$wsAB.Range("B81").Value = "contacts{}.city"
$wsAB.Range("C81").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsAB.Range("A82").Value = "contact_administrative_area"
$wsAB.Range("B82").Value = "contacts{}.administrative_area"
$wsAB.Range("C82").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsAB.Range("A83").Value = "contact_postalcode"
$wsAB.Range("B83").Value = "contacts{}.postalcode"
$wsAB.Range("C83").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsAB.Range("A84").Value = "contact_country"
$wsAB.Range("B84").Value = "contacts{}.country"
$wsAB.Range("C84").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

# Insert 3 rows above the contact block (pushes it to 81-87) and fill
# them in with the online_resource_description fields.
$wsAB.Rows("78:80").Insert()

$wsAB.Range("A78").Value = "online_resource_description"
$wsAB.Range("B78").Value = "online_resource_description_en"
$wsAB.Range("C78").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsAB.Range("A79").Value = "online_resource_description_other_lang_locale"
$wsAB.Range("B79").Value = "online_resource_description_locale"
$wsAB.Range("C79").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsAB.Range("A80").Value = "online_resource_description_other_lang"
$wsAB.Range("B80").Value = "online_resource_description_fr"
$wsAB.Range("C80").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

# Leave the view scrolled/selected over the freshly-added block.
$wsAB.Range("A76:XFD87").Select() | Out-Null

# ---------------------------------------------------------------------
# Config_BC (sheet3)
# ---------------------------------------------------------------------
$wsBC = $wb.Worksheets.Item("Config_BC")
$wsBC.Activate()

$wsBC.Range("A69").Value = "online_resource_link"
$wsBC.Range("B69").Value = "online_resource_link"
$wsBC.Range("C69").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A70").Value = "online_resource_protocol"
$wsBC.Range("B70").Value = "online_resource_protocol"
$wsBC.Range("C70").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A71").Value = "online_resource_description"
$wsBC.Range("B71").Value = "online_resource_description_en"
$wsBC.Range("C71").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A72").Value = "online_resource_description_other_lang_locale"
$wsBC.Range("B72").Value = "online_resource_description_locale"
$wsBC.Range("C72").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A73").Value = "online_resource_description_other_lang"
$wsBC.Range("B73").Value = "online_resource_description_fr"
$wsBC.Range("C73").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A74").Value = "contact_phone"
$wsBC.Range("B74").Value = "contacts{}.phone"
$wsBC.Range("C74").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A75").Value = "contact_facsimile"
$wsBC.Range("B75").Value = "contacts{}.facsimile"
$wsBC.Range("C75").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A76").Value = "contact_address"
$wsBC.Range("B76").Value = "contacts{}.address"
$wsBC.Range("C76").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A77").Value = "contact_city"
$wsBC.Range("B77").Value = "contacts{}.city"
$wsBC.Range("C77").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A78").Value = "contact_administrative_area"
$wsBC.Range("B78").Value = "contacts{}.administrative_area"
$wsBC.Range("C78").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A79").Value = "contact_postalcode"
$wsBC.Range("B79").Value = "contacts{}.postalcode"
$wsBC.Range("C79").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

$wsBC.Range("A80").Value = "contact_country"
$wsBC.Range("B80").Value = "contacts{}.country"
$wsBC.Range("C80").Value = "GMD_CONTACT;GMD_CITEDRESPONSIBLEPARTY;GMD_DISTRIBUTOR"

# Final selection/view state - Config_BC ends up the active sheet/tab.
$wsBC.Range("C70").Select() | Out-Null
